# Rename the "Collection_SU" tab to "CRF_SU" (CRF specialization for SU).
# Renaming the sheet also auto-updates the sheet-scoped defined name
# (_xlnm._FilterDatabase) that references it, e.g.
#   Collection_SU!$A$1:$AK$238  ->  CRF_SU!$A$1:$AK$238
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_SU")
$ws.Name = "CRF_SU"
